$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Update header labels (B1/C1) on both sheets
$wsMeans.Range("B1").Value = "Rural Areas (National Average)"
$wsMeans.Range("C1").Value = "Rural Areas (State Average)"
$wsSD.Range("B1").Value = "Rural Areas (National Average) SD"
$wsSD.Range("C1").Value = "Rural Areas (State Average) SD"

# --- Means sheet updates ---
$meansB = @(84, 7.6, 8.2, 10, 67, 6.8, 5.1, 26, 0.32)
$meansC = @(94, 2, 3.9, 3.5, 66, 5.4, 3.9, 21, 0.27)

for ($i = 0; $i -lt $meansB.Length; $i++) {
    $row = $i + 2
    $wsMeans.Cells.Item($row, 2).Value = $meansB[$i]
    $wsMeans.Cells.Item($row, 3).Value = $meansC[$i]
    $wsMeans.Cells.Item($row, 4).Value = "#NUM!"
    $wsMeans.Cells.Item($row, 5).Value = "#NUM!"
    $wsMeans.Cells.Item($row, 6).Value = "#NUM!"
    $wsMeans.Cells.Item($row, 7).Value = "#NUM!"
}

# --- Standard Deviations sheet updates ---
$sdB = @(19, 16, 12, 15, 28, 7.6, 6.4, 8.6, 0.14)
$sdC = @(6.9, 4.3, 4.6, 5.5, 19, 5.9, 4.6, 2.9, 0.051)

for ($i = 0; $i -lt $sdB.Length; $i++) {
    $row = $i + 2
    $wsSD.Cells.Item($row, 2).Value = $sdB[$i]
    $wsSD.Cells.Item($row, 3).Value = $sdC[$i]
    $wsSD.Cells.Item($row, 4).Value = 0
    $wsSD.Cells.Item($row, 5).Value = 0
    $wsSD.Cells.Item($row, 6).Value = 0
    $wsSD.Cells.Item($row, 7).Value = 0
}
